$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 23356
$ws.Range("I19").Value = 1533.3334
$ws.Range("J19").Value = 34267.332
$ws.Range("K19").Value = 1533.3334
$ws.Range("L19").Value = 34267.332
$ws.Range("M19").Value = -1358.3334
$ws.Range("N19").Value = -34617.332

$ws.Range("H98").Value = 1660
$ws.Range("I98").Value = 300
$ws.Range("K98").Value = 300
$ws.Range("M98").Value = 1198

$ws.Range("H113").Value = 2201
$ws.Range("I113").Value = 2176.25
$ws.Range("J113").Value = 2300
$ws.Range("K113").Value = 2176.25
$ws.Range("L113").Value = 2300
$ws.Range("M113").Value = 1077.75
$ws.Range("N113").Value = -8808

$ws.Range("H122").Value = 1660
$ws.Range("I122").Value = 300
$ws.Range("K122").Value = 900
$ws.Range("M122").Value = 1550

$ws.Range("H138").Value = 2128.9553
$ws.Range("J138").Value = 2110.9575
$ws.Range("L138").Value = 6332.872499999999
$ws.Range("N138").Value = -16612.8725

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1813.5
$ws.Range("I2").Value = 1951.9166
$ws.Range("J2").Value = 983
$ws.Range("K2").Value = 1951.9166
$ws.Range("L2").Value = 983
$ws.Range("M2").Value = -1838.9166
$ws.Range("N2").Value = -1209

$ws.Range("H32").Value = 436283.88
$ws.Range("I32").Value = 515817.28
$ws.Range("J32").Value = 18733.584
$ws.Range("K32").Value = 515817.28
$ws.Range("L32").Value = 18733.584
$ws.Range("M32").Value = -515530.28
$ws.Range("N32").Value = -19307.584

$ws.Range("H45").Value = 4492.6665
$ws.Range("I45").Value = 4152
$ws.Range("J45").Value = 4833.3335
$ws.Range("K45").Value = 4152
$ws.Range("L45").Value = 4833.3335
$ws.Range("M45").Value = -3775
$ws.Range("N45").Value = -5587.3335

$ws.Range("H61").Value = 3832.4211
$ws.Range("I61").Value = 3491.6
$ws.Range("J61").Value = 4211.1113
$ws.Range("K61").Value = 3491.6
$ws.Range("L61").Value = 4211.1113
$ws.Range("M61").Value = -3279.6
$ws.Range("N61").Value = -4635.1113

$ws.Range("H74").Value = 1961.2424
$ws.Range("I74").Value = 1627.2963
$ws.Range("J74").Value = 3464
$ws.Range("K74").Value = 1627.2963
$ws.Range("L74").Value = 3464
$ws.Range("M74").Value = -753.2963
$ws.Range("N74").Value = -5212

$ws.Range("H77").Value = 1961.2424
$ws.Range("I77").Value = 1627.2963
$ws.Range("J77").Value = 3464
$ws.Range("K77").Value = 8136.4815
$ws.Range("L77").Value = 17320
$ws.Range("M77").Value = -3768.4815
$ws.Range("N77").Value = -26056

$ws.Range("H97").Value = 1175
$ws.Range("I97").Value = 1175
$ws.Range("K97").Value = 1175
$ws.Range("M97").Value = -679

$ws.Range("H116").Value = 1813.5
$ws.Range("I116").Value = 1951.9166
$ws.Range("J116").Value = 983
$ws.Range("K116").Value = 1951.9166
$ws.Range("L116").Value = 983
$ws.Range("M116").Value = 342.0834
$ws.Range("N116").Value = -5571

$ws.Range("H136").Value = 3832.4211
$ws.Range("I136").Value = 3491.6
$ws.Range("J136").Value = 4211.1113
$ws.Range("K136").Value = 10474.8
$ws.Range("L136").Value = 12633.3339
$ws.Range("M136").Value = -7924.799999999999
$ws.Range("N136").Value = -17733.3339

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1813.5
$ws.Range("I3").Value = 1951.9166
$ws.Range("J3").Value = 983
$ws.Range("K3").Value = 1951.9166
$ws.Range("L3").Value = 983
$ws.Range("M3").Value = -1837.9166
$ws.Range("N3").Value = -1211

$ws.Range("H81").Value = 55483.223
$ws.Range("J81").Value = 55483.223
$ws.Range("L81").Value = 55483.223
$ws.Range("N81").Value = -57605.223

$ws.Range("H84").Value = 55483.223
$ws.Range("J84").Value = 55483.223
$ws.Range("L84").Value = 166449.669
$ws.Range("N84").Value = -177057.669

$ws.Range("H99").Value = 2330
$ws.Range("I99").Value = 2510
$ws.Range("K99").Value = 2510
$ws.Range("M99").Value = -1012

$ws.Range("H107").Value = 144564.28
$ws.Range("I107").Value = 168158.33
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 168158.33
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = -166238.33
$ws.Range("N107").Value = -6840

$ws.Range("H134").Value = 2900.348
$ws.Range("I134").Value = 2710.9473
$ws.Range("J134").Value = 3800
$ws.Range("K134").Value = 8132.841899999999
$ws.Range("L134").Value = 11400
$ws.Range("M134").Value = -5597.841899999999
$ws.Range("N134").Value = -16470

$ws.Range("H135").Value = 38608.43
$ws.Range("J135").Value = 38608.43
$ws.Range("L135").Value = 38608.43
$ws.Range("N135").Value = -48748.43

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2180.2727
$ws.Range("I16").Value = 1707.3334
$ws.Range("K16").Value = 1707.3334
$ws.Range("M16").Value = -1420.3334

$ws.Range("H31").Value = 5659.647
$ws.Range("I31").Value = 1296.6
$ws.Range("J31").Value = 11892.571
$ws.Range("K31").Value = 1296.6
$ws.Range("L31").Value = 11892.571
$ws.Range("M31").Value = -1001.6
$ws.Range("N31").Value = -12482.571

$ws.Range("H34").Value = 5659.647
$ws.Range("I34").Value = 1296.6
$ws.Range("J34").Value = 11892.571
$ws.Range("K34").Value = 1296.6
$ws.Range("L34").Value = 11892.571
$ws.Range("M34").Value = -1094.6
$ws.Range("N34").Value = -12296.571

$ws.Range("H58").Value = 1490.25
$ws.Range("I58").Value = 1165.625
$ws.Range("J58").Value = 2139.5
$ws.Range("K58").Value = 1165.625
$ws.Range("L58").Value = 2139.5
$ws.Range("M58").Value = -962.625
$ws.Range("N58").Value = -2545.5

$ws.Range("H113").Value = 2180.2727
$ws.Range("I113").Value = 1707.3334
$ws.Range("K113").Value = 1707.3334
$ws.Range("M113").Value = 462.6666

$ws.Range("H136").Value = 1490.25
$ws.Range("I136").Value = 1165.625
$ws.Range("J136").Value = 2139.5
$ws.Range("K136").Value = 3496.875
$ws.Range("L136").Value = 6418.5
$ws.Range("M136").Value = -946.875
$ws.Range("N136").Value = -11518.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 235.92857
$ws.Range("J12").Value = 258.16666
$ws.Range("L12").Value = 774.4999799999999
$ws.Range("N12").Value = -1120.49998

$ws.Range("H122").Value = 11792.556
$ws.Range("I122").Value = 422.5
$ws.Range("K122").Value = 3802.5
$ws.Range("M122").Value = -1352.5

$ws.Range("H131").Value = 839.4722
$ws.Range("I131").Value = 244
$ws.Range("J131").Value = 1068.5
$ws.Range("K131").Value = 732
$ws.Range("L131").Value = 3205.5
$ws.Range("M131").Value = 4308
$ws.Range("N131").Value = -13285.5

$ws.Range("H137").Value = 5562537.5
$ws.Range("J137").Value = 4356.737
$ws.Range("L137").Value = 13070.211
$ws.Range("N137").Value = -23270.211

$ws.Range("H140").Value = 1374.3182
$ws.Range("J140").Value = 2013.25
$ws.Range("L140").Value = 6039.75
$ws.Range("N140").Value = -16399.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 21000
$ws.Range("J42").Value = 21000
$ws.Range("L42").Value = 21000
$ws.Range("N42").Value = -21970

$ws.Range("H97").Value = 3176.6667
$ws.Range("I97").Value = 3072.375
$ws.Range("J97").Value = 4011
$ws.Range("K97").Value = 3072.375
$ws.Range("L97").Value = 4011
$ws.Range("M97").Value = -2576.375
$ws.Range("N97").Value = -5003

$ws.Range("H115").Value = 21000
$ws.Range("J115").Value = 21000
$ws.Range("L115").Value = 21000
$ws.Range("N115").Value = -23350

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2549.7778
$ws.Range("I61").Value = 2173.5
$ws.Range("J61").Value = 5560
$ws.Range("K61").Value = 2173.5
$ws.Range("L61").Value = 5560
$ws.Range("M61").Value = -1971.5
$ws.Range("N61").Value = -5964

$ws.Range("H100").Value = 2279.5
$ws.Range("I100").Value = 1781.2
$ws.Range("J100").Value = 2777.8
$ws.Range("K100").Value = 1781.2
$ws.Range("L100").Value = 2777.8
$ws.Range("M100").Value = -1240.2
$ws.Range("N100").Value = -3859.8

$ws.Range("H113").Value = 2549.7778
$ws.Range("I113").Value = 2173.5
$ws.Range("J113").Value = 5560
$ws.Range("K113").Value = 2173.5
$ws.Range("L113").Value = 5560
$ws.Range("M113").Value = -3.5
$ws.Range("N113").Value = -9900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2417.04
$ws.Range("I136").Value = 2142.75
$ws.Range("J136").Value = 9000
$ws.Range("K136").Value = 6428.25
$ws.Range("L136").Value = 27000
$ws.Range("M136").Value = -3878.25
$ws.Range("N136").Value = -32100
